$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 83.85416666666666
$ws.Range("C2").Value = 81.25
$ws.Range("D2").Value = 58.07291666666667
$ws.Range("E2").Value = 54.6875
$ws.Range("F2").Value = 53.125
$ws.Range("G2").Value = 52.08333333333333
$ws.Range("H2").Value = 52.08333333333333
$ws.Range("I2").Value = 51.82291666666667
$ws.Range("J2").Value = 50.52083333333333
$ws.Range("K2").Value = 49.21875
$ws.Range("L2").Value = 48.69791666666666
$ws.Range("M2").Value = 48.17708333333334
$ws.Range("N2").Value = 48.17708333333334
$ws.Range("O2").Value = 48.17708333333334
$ws.Range("P2").Value = 48.17708333333334
$ws.Range("Q2").Value = 48.17708333333334
$ws.Range("R2").Value = 48.17708333333334
$ws.Range("S2").Value = 48.17708333333334
$ws.Range("T2").Value = 48.17708333333334
$ws.Range("U2").Value = 48.17708333333334
